$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.512.09'
$ws.Range('E2').Value = '  +0.07%  '
$ws.Range('D3').Value = '2.080.45'
$ws.Range('E3').Value = '  +0.15%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = '234.53'
$ws.Range('E5').Value = '  -0.27%  '
$ws.Range('D6').Value = '0.625'
$ws.Range('E6').Value = '  +1.25%  '
$ws.Range('E7').Value = '  -0.03%  '
$ws.Range('E8').Value = '  -1.08%  '
$ws.Range('D9').Value = '0.388'
$ws.Range('E9').Value = '  +0.74%  '
$ws.Range('D10').Value = '0.0778'
$ws.Range('E10').Value = '  +1.89%  '
$ws.Range('E11').Value = '  +1.13%  '
$ws.Range('D12').Value = '2.387.75'
$ws.Range('E12').Value = '  +0.13%  '
$ws.Range('E13').Value = '  -1.19%  '
$ws.Range('D14').Value = '20.78'
$ws.Range('E14').Value = '  -1.32%  '
$ws.Range('E15').Value = '  +0.36%  '
$ws.Range('E16').Value = '  +0.36%  '
$ws.Range('D17').Value = '2.088.89'
$ws.Range('E17').Value = '  -0.68%  '
$ws.Range('D18').Value = '37.455.91'
$ws.Range('E18').Value = '  -0.70%  '
$ws.Range('D19').Value = '6.17'
$ws.Range('E19').Value = '  -1.34%  '
$ws.Range('D20').Value = '69.63'
$ws.Range('E20').Value = '  -0.91%  '
$ws.Range('E21').Value = '  +0.31%  '
$ws.Range('D22').Value = '226.44'
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('E24').Value = '  +2.62%  '
$ws.Range('E25').Value = '  -2.12%  '
$ws.Range('D26').Value = '168.36'
$ws.Range('E26').Value = '  +1.06%  '
$ws.Range('D27').Value = '8.88'
$ws.Range('E27').Value = '  -0.73%  '
$ws.Range('D28').Value = '1.43'
$ws.Range('E28').Value = '  -4.64%  '
$ws.Range('E29').Value = '  +2.80%  '
$ws.Range('D30').Value = '19.17'
$ws.Range('E30').Value = '  -0.62%  '
$ws.Range('E31').Value = '  -0.39%  '
$ws.Range('D32').Value = '4.62'
$ws.Range('E32').Value = '  +2.28%  '
$ws.Range('E33').Value = '  -1.12%  '
$ws.Range('E34').Value = '  +0.63%  '
$ws.Range('D35').Value = '2.54'
$ws.Range('E35').Value = '  -1.50%  '
$ws.Range('D36').Value = '3.45'
$ws.Range('E36').Value = '  +2.85%  '
$ws.Range('E37').Value = '  +0.54%  '
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('D39').Value = '5.57'
$ws.Range('E39').Value = '  -5.25%  '
$ws.Range('E40').Value = '  -0.35%  '
$ws.Range('B41').Value = 'Cronos'
$ws.Range('C41').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D41').Value = '0.0960'
$ws.Range('E41').Value = '  +1.07%  '
$ws.Range('B42').Value = 'Maker'
$ws.Range('C42').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D42').Value = '1.489.02'
$ws.Range('E42').Value = '  +2.18%  '
$ws.Range('E43').Value = '  +1.18%  '
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('E45').Value = '  -2.25%  '
$ws.Range('D46').Value = '4.20'
$ws.Range('E46').Value = '  -9.85%  '
$ws.Range('E47').Value = '  +0.61%  '
$ws.Range('D48').Value = '15.52'
$ws.Range('E48').Value = '  -2.27%  '
$ws.Range('E49').Value = '  +0.15%  '
$ws.Range('E50').Value = '  +1.37%  '
$ws.Range('D51').Value = '2.274.66'
$ws.Range('E51').Value = '  +0.17%  '
